$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.466.77"
$ws.Range("E2").Value = "'  -2.86%  "

# Row 3
$ws.Range("D3").Value = "'1.805.73"
$ws.Range("E3").Value = "'  -2.59%  "

# Row 4
$ws.Range("E4").Value = "'  +0.76%  "

# Row 5
$ws.Range("E5").Value = "'  +0.72%  "

# Row 6
$ws.Range("D6").Value = "'308.67"

# Row 7
$ws.Range("D7").Value = "'0.4541"
$ws.Range("E7").Value = "'  -1.45%  "

# Row 8
$ws.Range("D8").Value = "'0.3658"
$ws.Range("E8").Value = "'  -1.29%  "

# Row 9
$ws.Range("D9").Value = "'0.07124"
$ws.Range("E9").Value = "'  -2.33%  "

# Row 10
$ws.Range("D10").Value = "'0.8749"
$ws.Range("E10").Value = "'  -1.25%  "

# Row 11
$ws.Range("D11").Value = "'0.07745"
$ws.Range("E11").Value = "'  -0.79%  "

# Row 12
$ws.Range("D12").Value = "'19.35"
$ws.Range("E12").Value = "'  -3.50%  "

# Row 13
$ws.Range("D13").Value = "'1.832.95"
$ws.Range("E13").Value = "'  -2.54%  "

# Row 14
$ws.Range("D14").Value = "'5.267"
$ws.Range("E14").Value = "'  -2.19%  "

# Row 15
$ws.Range("D15").Value = "'6.343"
$ws.Range("E15").Value = "'  -2.76%  "

# Row 16
$ws.Range("D16").Value = "'86.13"
$ws.Range("E16").Value = "'  -5.80%  "

# Row 17
$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "'  +0.86%  "

# Row 18
$ws.Range("D18").Value = "'0.000008581"
$ws.Range("E18").Value = "'  -3.89%  "

# Row 19
$ws.Range("E19").Value = "'  +0.52%  "

# Row 20
$ws.Range("D20").Value = "'26.516.07"
$ws.Range("E20").Value = "'  -2.70%  "

# Row 21
$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "'  -3.54%  "

# Row 22
$ws.Range("D22").Value = "'4.973"
$ws.Range("E22").Value = "'  -2.59%  "

# Row 24
$ws.Range("D24").Value = "'1.977"
$ws.Range("E24").Value = "'  +2.70%  "

# Row 25
$ws.Range("D25").Value = "'150.98"
$ws.Range("E25").Value = "'  -0.29%  "

# Row 26
$ws.Range("D26").Value = "'17.92"
$ws.Range("E26").Value = "'  -2.54%  "

# Row 27
$ws.Range("D27").Value = "'2.002"
$ws.Range("E27").Value = "'  -2.65%  "

# Row 28
$ws.Range("D28").Value = "'112.86"
$ws.Range("E28").Value = "'  -2.61%  "

# Row 29
$ws.Range("D29").Value = "'4.847"
$ws.Range("E29").Value = "'  -4.20%  "

# Row 30
$ws.Range("D30").Value = "'0.08664"
$ws.Range("E30").Value = "'  -1.58%  "

# Row 31
$ws.Range("D31").Value = "'3.051"
$ws.Range("E31").Value = "'  -1.28%  "

# Row 32
$ws.Range("B32").Value = "'ImmutableX"
$ws.Range("C32").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7277"
$ws.Range("E32").Value = "'  -5.21%  "

# Row 33
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.439"
$ws.Range("E33").Value = "'  -1.33%  "

# Row 34
$ws.Range("D34").Value = "'1.110"
$ws.Range("E34").Value = "'  -5.15%  "

# Row 35
$ws.Range("D35").Value = "'1.008"
$ws.Range("E35").Value = "'  +0.79%  "

# Row 36
$ws.Range("D36").Value = "'2.554"
$ws.Range("E36").Value = "'  -6.52%  "

# Row 37
$ws.Range("D37").Value = "'1.081"
$ws.Range("E37").Value = "'  -0.05%  "

# Row 38
$ws.Range("D38").Value = "'0.01927"
$ws.Range("E38").Value = "'  -1.09%  "

# Row 39
$ws.Range("D39").Value = "'0.05095"
$ws.Range("E39").Value = "'  -2.99%  "

# Row 40
$ws.Range("D40").Value = "'2.879"
$ws.Range("E40").Value = "'  -2.34%  "

# Row 41
$ws.Range("D41").Value = "'6.954"
$ws.Range("E41").Value = "'  -1.57%  "

# Row 42
$ws.Range("D42").Value = "'0.4967"
$ws.Range("E42").Value = "'  -2.71%  "

# Row 43
$ws.Range("D43").Value = "'0.1565"
$ws.Range("E43").Value = "'  -3.93%  "

# Row 44
$ws.Range("D44").Value = "'8.115"
$ws.Range("E44").Value = "'  -3.20%  "

# Row 45
$ws.Range("E45").Value = "'  +0.82%  "

# Row 46
$ws.Range("D46").Value = "'0.4615"
$ws.Range("E46").Value = "'  -3.61%  "

# Row 47
$ws.Range("B47").Value = "'Quant"
$ws.Range("C47").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'102.19"
$ws.Range("E47").Value = "'  +0.03%  "

# Row 48
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.937"
$ws.Range("E48").Value = "'  -3.86%  "

# Row 49
$ws.Range("D49").Value = "'1.588"
$ws.Range("E49").Value = "'  -3.26%  "

# Row 50
$ws.Range("D50").Value = "'0.06000"
$ws.Range("E50").Value = "'  -3.45%  "

# Row 51
$ws.Range("D51").Value = "'63.78"
$ws.Range("E51").Value = "'  -2.79%  "
